$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.013793926068955
$ws.Range("D2").Value = 1.020048264153338
$ws.Range("E2").Value = 1.015670437013835
$ws.Range("I2").Value = 1.025677096293251
$ws.Range("J2").Value = 1.019028854079417
$ws.Range("K2").Value = 1.022889703837036
$ws.Range("L2").Value = 1.018524900757735
$ws.Range("N2").Value = 1.010482264980734

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.014560158083
$ws.Range("D3").Value = 1.020571571607653
$ws.Range("E3").Value = 1.01631489921221
$ws.Range("I3").Value = 1.025731484584719
$ws.Range("J3").Value = 1.01943035070711
$ws.Range("K3").Value = 1.023219944023796
$ws.Range("L3").Value = 1.018975022064269
$ws.Range("N3").Value = 1.010615547099763

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.015056621355488
$ws.Range("D4").Value = 1.0209104740811
$ws.Range("E4").Value = 1.016732888917631
$ws.Range("I4").Value = 1.025765382872579
$ws.Range("J4").Value = 1.01969015008195
$ws.Range("K4").Value = 1.023433207292328
$ws.Range("L4").Value = 1.019266574955233
$ws.Range("N4").Value = 1.010701766359047

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.015265490542645
$ws.Range("D5").Value = 1.021053015495152
$ws.Range("E5").Value = 1.016908844252872
$ws.Range("I5").Value = 1.025779323072867
$ws.Range("J5").Value = 1.019799369339803
$ws.Range("K5").Value = 1.023522760176482
$ws.Range("L5").Value = 1.019389212717889
$ws.Range("N5").Value = 1.010738006891468

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.015300569690796
$ws.Range("D6").Value = 1.021076952684722
$ws.Range("E6").Value = 1.016938401500015
$ws.Range("I6").Value = 1.025781645456922
$ws.Range("J6").Value = 1.019817707669409
$ws.Range("K6").Value = 1.02353779041835
$ws.Range("L6").Value = 1.019409808135243
$ws.Range("N6").Value = 1.010744091467512

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.015059411665167
$ws.Range("D7").Value = 1.02091237846648
$ws.Range("E7").Value = 1.016735239130595
$ws.Range("I7").Value = 1.025765570363717
$ws.Range("J7").Value = 1.01969160947836
$ws.Range("K7").Value = 1.023434404308855
$ws.Range("L7").Value = 1.019268213378111
$ws.Range("N7").Value = 1.010702250631187

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.014052739945947
$ws.Range("D8").Value = 1.020225057599428
$ws.Range("E8").Value = 1.015888031974956
$ws.Range("I8").Value = 1.025695744690433
$ws.Range("J8").Value = 1.019164540027655
$ws.Range("K8").Value = 1.023001397028323
$ws.Range("L8").Value = 1.018676959503296
$ws.Range("N8").Value = 1.01052731280066

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.012283996683516
$ws.Range("D9").Value = 1.01901620227694
$ws.Range("E9").Value = 1.014402737519766
$ws.Range("I9").Value = 1.025562819711759
$ws.Range("J9").Value = 1.018235870720592
$ws.Range("K9").Value = 1.022235197315357
$ws.Range("L9").Value = 1.017637416281911
$ws.Range("N9").Value = 1.010218894241586

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.011108408088345
$ws.Range("D10").Value = 1.01821196073472
$ws.Range("E10").Value = 1.013417768607924
$ws.Range("I10").Value = 1.025467600937058
$ws.Range("J10").Value = 1.017616905359605
$ws.Range("K10").Value = 1.021722341455993
$ws.Range("L10").Value = 1.016946039847802
$ws.Range("N10").Value = 1.010013206632241

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.010600237626655
$ws.Range("D11").Value = 1.017864135158167
$ws.Range("E11").Value = 1.012992532047668
$ws.Range("I11").Value = 1.02542481326569
$ws.Range("J11").Value = 1.017348939311283
$ws.Range("K11").Value = 1.021499799700126
$ws.Range("L11").Value = 1.016647077525566
$ws.Range("N11").Value = 1.009924130024823

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01041161277897
$ws.Range("D12").Value = 1.017735001892168
$ws.Range("E12").Value = 1.012834771903383
$ws.Range("I12").Value = 1.025408686767084
$ws.Range("J12").Value = 1.017249413750399
$ws.Range("K12").Value = 1.021417068371979
$ws.Range("L12").Value = 1.016536092518
$ws.Range("N12").Value = 1.009891041632509

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.010452067430553
$ws.Range("D13").Value = 1.017762698457446
$ws.Range("E13").Value = 1.01286860328333
$ws.Range("I13").Value = 1.025412156499441
$ws.Range("J13").Value = 1.017270761892427
$ws.Range("K13").Value = 1.021434817640167
$ws.Range("L13").Value = 1.016559896306762
$ws.Range("N13").Value = 1.009898139259054

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.010584643127656
$ws.Range("D14").Value = 1.0178534596314
$ws.Range("E14").Value = 1.0129794876231
$ws.Range("I14").Value = 1.02542348499842
$ws.Range("J14").Value = 1.017340712306243
$ws.Range("K14").Value = 1.021492962510221
$ws.Range("L14").Value = 1.01663790217776
$ws.Range("N14").Value = 1.009921394953103

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.010666344979696
$ws.Range("D15").Value = 1.017909389214088
$ws.Range("E15").Value = 1.013047832583145
$ws.Range("I15").Value = 1.025430433973751
$ws.Range("J15").Value = 1.017383812295339
$ws.Range("K15").Value = 1.021528778329882
$ws.Range("L15").Value = 1.016685972552733
$ws.Range("N15").Value = 1.009935723386767

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.011142152129631
$ws.Range("D16").Value = 1.018235053722858
$ws.Range("E16").Value = 1.013446016917354
$ws.Range("I16").Value = 1.025470407864246
$ws.Range("J16").Value = 1.017634690561183
$ws.Range("K16").Value = 1.021737101004919
$ws.Range("L16").Value = 1.016965889728307
$ws.Range("N16").Value = 1.010019118132853

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.011440846963881
$ws.Range("D17").Value = 1.018439447405133
$ws.Range("E17").Value = 1.01369612666564
$ws.Range("I17").Value = 1.025495065998497
$ws.Range("J17").Value = 1.01779207428526
$ws.Range("K17").Value = 1.021867651043158
$ws.Range("L17").Value = 1.017141584788225
$ws.Range("N17").Value = 1.010071426460765

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.011615154140687
$ws.Range("D18").Value = 1.018558706810633
$ws.Range("E18").Value = 1.013842133022439
$ws.Range("I18").Value = 1.02550929836483
$ws.Range("J18").Value = 1.017883878280057
$ws.Range("K18").Value = 1.021943753053054
$ws.Range("L18").Value = 1.017244104027421
$ws.Range("N18").Value = 1.010101935782033

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.011674602491142
$ws.Range("D19").Value = 1.018599377911759
$ws.Range("E19").Value = 1.013891937972986
$ws.Range("I19").Value = 1.02551412571761
$ws.Range("J19").Value = 1.017915181828418
$ws.Range("K19").Value = 1.021969694089316
$ws.Range("L19").Value = 1.017279067053739
$ws.Range("N19").Value = 1.010112338446175

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.011408791180248
$ws.Range("D20").Value = 1.018417513741933
$ws.Range("E20").Value = 1.013669279666561
$ws.Range("I20").Value = 1.025492435955763
$ws.Range("J20").Value = 1.017775188001345
$ws.Range("K20").Value = 1.021853648967525
$ws.Range("L20").Value = 1.017122730296427
$ws.Range("N20").Value = 1.010065814396976

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.0105455992213
$ws.Range("D21").Value = 1.017826730926677
$ws.Range("E21").Value = 1.012946829641299
$ws.Range("I21").Value = 1.025420155470982
$ws.Range("J21").Value = 1.017320113393908
$ws.Range("K21").Value = 1.021475842199527
$ws.Range("L21").Value = 1.016614929643717
$ws.Range("N21").Value = 1.009914546764714

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.010003642548644
$ws.Range("D22").Value = 1.017455657341106
$ws.Range("E22").Value = 1.012493706499807
$ws.Range("I22").Value = 1.025373360689526
$ws.Range("J22").Value = 1.017034042482262
$ws.Range("K22").Value = 1.021237899221228
$ws.Range("L22").Value = 1.016296019976106
$ws.Range("N22").Value = 1.009819431099563

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.010290870721381
$ws.Range("D23").Value = 1.017652334225293
$ws.Range("E23").Value = 1.012733809740103
$ws.Range("I23").Value = 1.025398295135386
$ws.Range("J23").Value = 1.017185688622326
$ws.Range("K23").Value = 1.021364074804179
$ws.Range("L23").Value = 1.016465044924327
$ws.Range("N23").Value = 1.009869854274592

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.011423275551202
$ws.Range("D24").Value = 1.018427424496807
$ws.Range("E24").Value = 1.01368141029448
$ws.Range("I24").Value = 1.025493624823379
$ws.Range("J24").Value = 1.017782818171756
$ws.Range("K24").Value = 1.021859976043938
$ws.Range("L24").Value = 1.01713124970852
$ws.Range("N24").Value = 1.010068350250936

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.012740637921531
$ws.Range("D25").Value = 1.019328436256194
$ws.Range("E25").Value = 1.014785809320876
$ws.Range("I25").Value = 1.025598350806098
$ws.Range("J25").Value = 1.018475934140528
$ws.Range("K25").Value = 1.022433646965284
$ws.Range("L25").Value = 1.017905878561795
$ws.Range("N25").Value = 1.010298643247282
